# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 44 -> 45 (ColumnWidth value chosen so the exported
#     OOXML <col width> rounds to exactly 45) ---
$ws.Columns("A").ColumnWidth = 44.14

# --- "Bad Drivers" block (row 3 = the single bad driver, row 4 = totals) ---
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.40.0.4"
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 282
$ws.Range("D3").Value = 98.90000000000001

$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 282

# --- "Good Drivers" block ---
# Row 12
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B12").Value = 11140
$ws.Range("D12").Value = 100

# E12 is a plain text cell like "2022-08-29" (General format, not a real
# date). Writing the literal text directly gets auto-coerced into a date
# serial by this engine, so instead build it as a formula and then
# Copy / PasteSpecial(values) over itself - that bakes in the literal text
# without the date auto-detection, keeping the original General/right-align
# style.
$ws.Range("E12").Formula = "=""2022-08-29"""
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B13").Value = 14487
$ws.Range("D13").Value = 100

$ws.Range("E13").Formula = "=""2022-05-23"""
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B14").Value = 265400
$ws.Range("D14").Value = 99.90000000000001

$ws.Range("E14").Formula = "=""2022-05-01"""
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

# Rows 15-17 previously held 3 more "good driver" entries; they are removed
# (dimension shrinks from J22 to J19) - fully clear them so no cell content
# or formatting remains.
$ws.Range("A15:J17").Clear()

$excel.CutCopyMode = $false
